$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the columns that get reshuffled across rows 2..16:
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $snapshot[$r] = @(
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 10).Value2,
        $ws.Cells.Item($r, 11).Value2,
        $ws.Cells.Item($r, 12).Value2,
        $ws.Cells.Item($r, 13).Value2,
        $ws.Cells.Item($r, 15).Value2,
        $ws.Cells.Item($r, 16).Value2
    )
}

# Map: destination row -> source row (which row's data it now receives)
$rowMap = @{
    2 = 16
    3 = 11
    4 = 15
    5 = 7
    6 = 3
    7 = 12
    8 = 6
    9 = 4
    10 = 5
    11 = 14
    12 = 10
    13 = 2
    14 = 13
    15 = 9
    16 = 8
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    $ws.Cells.Item($destRow, 4).Value = $vals[0]
    $ws.Cells.Item($destRow, 10).Value = $vals[1]
    $ws.Cells.Item($destRow, 11).Value = $vals[2]
    $ws.Cells.Item($destRow, 12).Value = $vals[3]
    $ws.Cells.Item($destRow, 13).Value = $vals[4]
    $ws.Cells.Item($destRow, 15).Value = $vals[5]
    $ws.Cells.Item($destRow, 16).Value = $vals[6]
}
